$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = -6.329399999999998
$ws.Range("A4").Value  = -21.13500000000001
$ws.Range("B4").Value  = 5.007000000000003
$ws.Range("D4").Value  = -6.926999999999996
$ws.Range("B5").Value  = 5.212600000000001
$ws.Range("A6").Value  = -21.29100000000001
$ws.Range("B6").Value  = 5.808099999999996
$ws.Range("A7").Value  = -21.26940000000001
$ws.Range("A8").Value  = -21.49230000000002
$ws.Range("B8").Value  = 4.940700000000001
$ws.Range("D9").Value  = -8.708399999999997
$ws.Range("D11").Value = -8.339499999999997
$ws.Range("D14").Value = -8.704599999999994
$ws.Range("A16").Value = -21.55530000000002
$ws.Range("B16").Value = 5.008400000000002
$ws.Range("D18").Value = -8.739899999999995
$ws.Range("A20").Value = -22.07020000000002
$ws.Range("A21").Value = -20.5711
$ws.Range("B22").Value = 5.263600000000004
$ws.Range("D25").Value = -7.474899999999995
